$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.541.31'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.911.75'
$ws.Range("E3").Value = '  +3.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.02'
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.95'
$ws.Range("E5").Value = '  +4.25%  '
$ws.Range("E6").Value = '  +5.62%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.24'
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("E9").Value = '  +5.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.26'
$ws.Range("E10").Value = '  +5.48%  '
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.101'
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.190.40'
$ws.Range("E13").Value = '  +3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.27'
$ws.Range("E14").Value = '  +7.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.699'
$ws.Range("E15").Value = '  +3.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.910.53'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.87'
$ws.Range("E17").Value = '  +3.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.576.16'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.33'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0826'
$ws.Range("E20").Value = '  +4.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '244.82'
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.70'
$ws.Range("E22").Value = '  +4.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.76'
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.17'
$ws.Range("E27").Value = '  +13.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.50'
$ws.Range("E28").Value = '  +7.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.47'
$ws.Range("E29").Value = '  +4.75%  '
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.963'
$ws.Range("E31").Value = '  +23.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.15'
$ws.Range("E32").Value = '  +3.96%  '
$ws.Range("E33").Value = '  +2.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.20'
$ws.Range("E34").Value = '  +4.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.02'
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.05'
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("E38").Value = '  +3.51%  '
$ws.Range("E39").Value = '  +3.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0212'
$ws.Range("E40").Value = '  +4.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '92.45'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0634'
$ws.Range("E42").Value = '  +14.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.80'
$ws.Range("E43").Value = '  +6.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.352.38'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  +4.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.80'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.31'
$ws.Range("E47").Value = '  +35.69%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.42'
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.58'
$ws.Range("E50").Value = '  +2.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.101.72'
$ws.Range("E51").Value = '  +3.33%  '
